$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.818.56'
$ws.Range("E2").Value = '  +0.87%  '
$ws.Range("D3").Value = '2.088.43'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '235.22'
$ws.Range("E5").Value = '  +0.15%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").Value = '58.93'
$ws.Range("E7").Value = '  +2.17%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '0.393'
$ws.Range("E9").Value = '  +1.16%  '
$ws.Range("E10").Value = '  +1.55%  '
$ws.Range("E11").Value = '  +2.61%  '
$ws.Range("D12").Value = '2.395.63'
$ws.Range("E12").Value = '  +0.28%  '
$ws.Range("D13").Value = '14.85'
$ws.Range("E13").Value = '  +2.88%  '
$ws.Range("D14").Value = '21.28'
$ws.Range("E14").Value = '  +2.04%  '
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("E16").Value = '  +2.16%  '
$ws.Range("D17").Value = '2.079.53'
$ws.Range("E17").Value = '  -0.20%  '
$ws.Range("D18").Value = '37.695.94'
$ws.Range("E18").Value = '  +0.64%  '
$ws.Range("D19").Value = '6.18'
$ws.Range("E19").Value = '  -0.12%  '
$ws.Range("D20").Value = '71.51'
$ws.Range("E20").Value = '  +2.62%  '
$ws.Range("D21").Value = '0.0₃0837'
$ws.Range("E21").Value = '  +1.94%  '
$ws.Range("D22").Value = '229.05'
$ws.Range("E22").Value = '  +1.12%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("E24").Value = '  -1.39%  '
$ws.Range("D25").Value = '2.43'
$ws.Range("E25").Value = '  +1.07%  '
$ws.Range("D26").Value = '169.94'
$ws.Range("E26").Value = '  +0.70%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.140'
$ws.Range("E27").Value = '  +6.43%  '
$ws.Range("D28").Value = '9.04'
$ws.Range("E28").Value = '  +1.49%  '
$ws.Range("D29").Value = '19.68'
$ws.Range("E29").Value = '  +2.65%  '
$ws.Range("E30").Value = '  -1.39%  '
$ws.Range("D32").Value = '4.72'
$ws.Range("E32").Value = '  +1.76%  '
$ws.Range("D33").Value = '0.0635'
$ws.Range("E33").Value = '  +2.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.70'
$ws.Range("E34").Value = '  +2.30%  '
$ws.Range("D35").Value = '2.51'
$ws.Range("E35").Value = '  -1.44%  '
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '1.84'
$ws.Range("E36").Value = '  +3.05%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '3.44'
$ws.Range("E37").Value = '  +1.34%  '
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("D39").Value = '5.42'
$ws.Range("E39").Value = '  -3.09%  '
$ws.Range("D40").Value = '0.0984'
$ws.Range("E40").Value = '  +3.05%  '
$ws.Range("D41").Value = '99.03'
$ws.Range("E41").Value = '  +1.86%  '
$ws.Range("E42").Value = '  -1.14%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '0.0216'
$ws.Range("E43").Value = '  +1.66%  '
$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").Value = '4.38'
$ws.Range("E44").Value = '  +4.36%  '
$ws.Range("D45").Value = '1.463.81'
$ws.Range("E45").Value = '  -1.75%  '
$ws.Range("E46").Value = '  +0.81%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = '1.07'
$ws.Range("E47").Value = '  +2.60%  '
$ws.Range("B48").Value = 'InjectiveProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '16.00'
$ws.Range("E48").Value = '  +3.02%  '
$ws.Range("E49").Value = '  +2.81%  '
$ws.Range("E50").Value = '  +2.00%  '
$ws.Range("D51").Value = '2.279.95'
$ws.Range("E51").Value = '  +0.23%  '
